# Update the numeric data range B2:H13 so it is stored as integer data
# (each value rounded to the nearest whole number), per commit:
# "test: Update references to integer format"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1
$lastCol = $used.Columns.Count + $used.Column - 1

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $dval = [double]$val
            if ($dval -ge 0) {
                $rounded = [Math]::Floor($dval + 0.5)
            } else {
                $rounded = [Math]::Ceiling($dval - 0.5)
            }
            $cell.Value2 = $rounded
        }
    }
}
